# AFDP-3458: Add new MyDocuments module
# - Add creator as default assignee for PERSONAL repositories (deny read access rule)
# - Deny read access to * participant for PERSONAL document repositories
#
# This inserts a new rule row right after the existing "Document Repository -
# default read access" rule (row 59) on Sheet1, pushing the rest of the
# DocumentRepository rules (old rows 60-67) down by one (new rows 61-68).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the existing DocumentRepository rule rows (60-67) down to (61-68),
# working from the bottom up so we never overwrite a source row before it
# has been copied. Only columns B:G carry data in this block (column A is
# blank/untouched throughout this range of rows).
for ($r = 67; $r -ge 60; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("B" + $srcRow + ":G" + $srcRow)
    $dst = $ws.Range("B" + $dstRow + ":G" + $dstRow)
    # Clear destination first - Copy() only overwrites cells that actually
    # have content in the source, so stale values would otherwise survive.
    $dst.ClearContents()
    $src.Copy($dst)
}

# Row 60 now still holds a duplicate of the old row 60 (now also at row 61).
# Clear it and build the new rule in its place, re-using the formatting
# already present on neighboring DocumentRepository rows so style indices
# line up (C/D/E/F use the "restricted flag" style, B/G use the plain
# wrap-text style).
$ws.Range("B60:G60").ClearContents()
$ws.Range("C65:F65").Copy($ws.Range("C60:F60"))
$ws.Range("B61").Copy($ws.Range("B60"))
$ws.Range("B61").Copy($ws.Range("G60"))

$ws.Range("B60").Value = "Document Repository -deny read access"
$ws.Range("C60").Value = "DOC_REPO"
$ws.Range("D60").Value = "repositoryType == 'PERSONAL'"
$ws.Range("G60").Value = "mandatory deny read to *"

# Restore row heights for the new/shifted rows.
$ws.Rows.Item(60).RowHeight = 30
$ws.Rows.Item(61).RowHeight = 45
$ws.Rows.Item(62).RowHeight = 30
$ws.Rows.Item(63).RowHeight = 30
$ws.Rows.Item(64).RowHeight = 30
$ws.Rows.Item(65).RowHeight = 30
$ws.Rows.Item(66).RowHeight = 45
$ws.Rows.Item(67).RowHeight = 45
$ws.Rows.Item(68).RowHeight = 60

# Match the saved selection/active cell from the edit.
$ws.Activate()
$ws.Range("I59").Select()

Write-Host "Inserted 'Document Repository -deny read access' rule at row 60"
